$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D (shifts D:K -> E:L), preserving formats
$ws.Columns("D").Insert()

# Copy number formats / styles from the (now shifted) neighboring column E into the
# freshly inserted column D, per contiguous data block (avoid touching the blank
# separator rows 36 and 78 which must remain absent from the sheet).
$ws.Range("E7:E35").Copy()
$ws.Range("D7:D35").PasteSpecial(-4122)
$ws.Range("E38:E77").Copy()
$ws.Range("D38:D77").PasteSpecial(-4122)
$ws.Range("E80:E102").Copy()
$ws.Range("D80:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Give the new column a sensible width similar to its neighbours
$ws.Columns("D").ColumnWidth = $ws.Columns("E").ColumnWidth

# Populate the new column D with the latest quarter's figures (new data),
# and restore the values for rows whose new-quarter value is unchanged
# from the previous quarter (these cells still need their value set because
# the Insert() operation leaves the new column's cells empty).
$ws.Range("D7").Value2 = 43373
$ws.Range("D8").Value2 = 8700
$ws.Range("D9").Value2 = 6900
$ws.Range("D10").Value2 = 1800
$ws.Range("D12").Value2 = "NA"
$ws.Range("D13").Value2 = 0
$ws.Range("D14").Value2 = 0
$ws.Range("D15").Value2 = 0
$ws.Range("D17").Value2 = 8000
$ws.Range("D18").Value2 = 700
$ws.Range("D20").Value2 = 0
$ws.Range("D21").Value2 = 800
$ws.Range("D22").Value2 = 0
$ws.Range("D23").Value2 = 700
$ws.Range("D24").Value2 = 200
$ws.Range("D25").Value2 = 0
$ws.Range("D26").Value2 = 500
$ws.Range("D27").Value2 = 500
$ws.Range("D28").Value2 = 0
$ws.Range("D29").Value2 = "NA"
$ws.Range("D30").Value2 = 0
$ws.Range("D31").Value2 = 0
$ws.Range("D32").Value2 = 0
$ws.Range("D33").Value2 = 500
$ws.Range("D34").Value2 = 0
$ws.Range("D35").Value2 = 500
$ws.Range("D38").Value2 = 43373
$ws.Range("D41").Value2 = 400
$ws.Range("D42").Value2 = 0
$ws.Range("D43").Value2 = 6800
$ws.Range("D44").Value2 = 12400
$ws.Range("D45").Value2 = 400
$ws.Range("D46").Value2 = 20100
$ws.Range("D47").Value2 = 0
$ws.Range("D48").Value2 = 4200
$ws.Range("D49").Value2 = 400
$ws.Range("D50").Value2 = 0
$ws.Range("D51").Value2 = 0
$ws.Range("D52").Value2 = 300
$ws.Range("D53").Value2 = 0
$ws.Range("D54").Value2 = 25100
$ws.Range("D57").Value2 = 1000
$ws.Range("D58").Value2 = 200
$ws.Range("D59").Value2 = 900
$ws.Range("D60").Value2 = 2100
$ws.Range("D61").Value2 = 0
$ws.Range("D62").Value2 = 100
$ws.Range("D63").Value2 = 0
$ws.Range("D64").Value2 = 0
$ws.Range("D65").Value2 = 0
$ws.Range("D66").Value2 = 2100
$ws.Range("D68").Value2 = 0
$ws.Range("D69").Value2 = 0
$ws.Range("D70").Value2 = 0
$ws.Range("D71").Value2 = 0
$ws.Range("D72").Value2 = 21800
$ws.Range("D73").Value2 = 0
$ws.Range("D74").Value2 = 0
$ws.Range("D75").Value2 = 0
$ws.Range("D76").Value2 = 22900
$ws.Range("D77").Value2 = 0
$ws.Range("D80").Value2 = 43373
$ws.Range("D81").Value2 = 500
$ws.Range("D83").Value2 = 100
$ws.Range("D84").Value2 = 0
$ws.Range("D85").Value2 = 0
$ws.Range("D86").Value2 = 0
$ws.Range("D87").Value2 = 0
$ws.Range("D88").Value2 = 0
$ws.Range("D89").Value2 = -1800
$ws.Range("D91").Value2 = 0
$ws.Range("D92").Value2 = 0
$ws.Range("D93").Value2 = 0
$ws.Range("D94").Value2 = 100
$ws.Range("D96").Value2 = 0
$ws.Range("D97").Value2 = 0
$ws.Range("D98").Value2 = 0
$ws.Range("D99").Value2 = 0
$ws.Range("D100").Value2 = -400
$ws.Range("D101").Value2 = 0
$ws.Range("D102").Value2 = -2100

$wb.Save()
